# InTheClearTests.xlsx edit:
#  - Insert a new "TB-2" bug-test row (pushing the old TB-2..TB-5 rows down,
#    renumbering the last of them to TB-6), adding 5 new shared strings.
#  - Fill in the previously-empty "Actual Result" cell for TE-6 (row 24).
#  - Update the active selection / scroll position to reflect where the
#    user ended up after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (pushes old rows 14-24 down to 15-25).
# The new row inherits formatting (wrap text, row style) from row 13 above it.
$ws.Rows("14:14").Insert()

# The test names in column A for the rows that got pushed down need to be
# renumbered sequentially again; the row that used to be TB-5 becomes TB-6.
$ws.Range("A18").Value = "TB-6"

# Populate the new TB-2 row with the new test case content.
$ws.Range("A14").Value = "TB-2"
$ws.Range("B14").Value = "Trip to a non-contiguous  IS State."
$ws.Range("C14").Value = "Trip should pop up an error if the route is longer than 51 hours"
$ws.Range("D14").Value = "Trip pops up an error if longer than 51 hours"

# This is a two-line wrapped row like the other "long" rows (height 34).
$ws.Range("A14:D14").RowHeight = 34

# The rest of the renumbered rows keep their original (shifted) text.
$ws.Range("A15").Value = "TB-3"
$ws.Range("A16").Value = "TB-4"
$ws.Range("A17").Value = "TB-5"

# Fill in the Actual Result for TE-6 (now row 24), which previously had no
# value in column D. Copy formatting from the neighboring cell first so the
# new cell matches the row's style, then set its text.
$ws.Range("C24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = "Route is given, map loads, and table still populates"

# Leave the selection where the user ended up after scrolling down to work
# on the newly added rows.
$ws.Range("E30").Select()
